$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Sachsen parliament period (119 -> 157), label + api url updated
$ws.Range("B3").Value = 157
$ws.Range("D3").Value = "Sachsen 2024 - 2029"
$ws.Range("E3").Value = "https://www.abgeordnetenwatch.de/api/v2/parliament-periods/157"

# Row 4: Brandenburg parliament period (154 -> 158), label + api url updated
$ws.Range("B4").Value = 158
$ws.Range("D4").Value = "Brandenburg 2024 - 2029"
$ws.Range("E4").Value = "https://www.abgeordnetenwatch.de/api/v2/parliament-periods/158"

# Row 5: Thüringen parliament period (121 -> 156), label + api url updated
$ws.Range("B5").Value = 156
$ws.Range("D5").Value = "Thüringen 2024 - 2029"
$ws.Range("E5").Value = "https://www.abgeordnetenwatch.de/api/v2/parliament-periods/156"
